# Comercializadora del Agro de Limarí - Nectarín: weekly refresh.
# Insert 7 new rows of data at the top of the data block (rows 58-64),
# pushing the existing rows (old 58-95) down to rows 65-102.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 blank rows before the current row 58 (shifts rows 58:95 -> 65:102).
$ws.Range("A58:A64").EntireRow.Insert()

# Common (constant) column values shared by every record row in this block.
$mercadoId   = 2
$mercado     = "Comercializadora del Agro de Limarí"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100103
$producto    = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria   = "Nectarín"

$rows = @(
    @{ Row=58; Fecha=44567; Variedad="Candy White";   Calidad="Primera";  Volumen=16; PMin=325000; PMax=330000; PProm=327500; Unidad="$/bins (420 kilos)";         Origen="Región de O'Higgins"; PKg=780;  KgUnidad=420 },
    @{ Row=59; Fecha=44567; Variedad="Candy White";   Calidad="Segunda";  Volumen=16; PMin=275000; PMax=280000; PProm=277500; Unidad="$/bins (420 kilos)";         Origen="Región de O'Higgins"; PKg=661;  KgUnidad=420 },
    @{ Row=60; Fecha=44567; Variedad="Early Diamond"; Calidad="Especial"; Volumen=10; PMin=415000; PMax=420000; PProm=417500; Unidad="$/bins (420 kilos)";         Origen="Región de O'Higgins"; PKg=994;  KgUnidad=420 },
    @{ Row=61; Fecha=44567; Variedad="Early Diamond"; Calidad="Primera";  Volumen=16; PMin=385000; PMax=390000; PProm=387500; Unidad="$/bins (420 kilos)";         Origen="Región de O'Higgins"; PKg=923;  KgUnidad=420 },
    @{ Row=62; Fecha=44567; Variedad="Early Diamond"; Calidad="Segunda";  Volumen=16; PMin=335000; PMax=340000; PProm=337500; Unidad="$/bins (420 kilos)";         Origen="Región de O'Higgins"; PKg=804;  KgUnidad=420 },
    @{ Row=63; Fecha=44567; Variedad="Magique";       Calidad="Especial"; Volumen=16; PMin=415000; PMax=420000; PProm=417500; Unidad="$/bins (420 kilos)";         Origen="Región de O'Higgins"; PKg=994;  KgUnidad=420 },
    @{ Row=64; Fecha=44567; Variedad="Magique";       Calidad="Primera";  Volumen=16; PMin=385000; PMax=390000; PProm=387500; Unidad="$/bins (420 kilos)";         Origen="Región de O'Higgins"; PKg=923;  KgUnidad=420 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Cells.Item($i, 1).Value  = $mercadoId
    $ws.Cells.Item($i, 2).Value  = $mercado
    $ws.Cells.Item($i, 3).Value  = $region
    $ws.Cells.Item($i, 4).Value  = $r.Fecha
    $ws.Cells.Item($i, 5).Value  = $codreg
    $ws.Cells.Item($i, 6).Value  = $tipo
    $ws.Cells.Item($i, 7).Value  = $productoId
    $ws.Cells.Item($i, 8).Value  = $producto
    $ws.Cells.Item($i, 9).Value  = $categoriaId
    $ws.Cells.Item($i, 10).Value = $categoria
    $ws.Cells.Item($i, 11).Value = $r.Variedad
    $ws.Cells.Item($i, 12).Value = $r.Calidad
    $ws.Cells.Item($i, 13).Value = $r.Volumen
    $ws.Cells.Item($i, 14).Value = $r.PMin
    $ws.Cells.Item($i, 15).Value = $r.PMax
    $ws.Cells.Item($i, 16).Value = $r.PProm
    $ws.Cells.Item($i, 17).Value = $r.Unidad
    $ws.Cells.Item($i, 18).Value = $r.Origen
    $ws.Cells.Item($i, 19).Value = $r.PKg
    $ws.Cells.Item($i, 20).Value = $r.KgUnidad
}
